$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.147.51"
$ws.Range("E2").Value = "  -3.44%  "
$ws.Range("D3").Value = "'1.603.38"
$ws.Range("E3").Value = "  -2.88%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'301.50"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("D7").Value = "'0.3781"
$ws.Range("E7").Value = "  -2.98%  "
$ws.Range("D8").Value = "'0.3651"
$ws.Range("E8").Value = "  -4.44%  "
$ws.Range("D9").Value = "'50.03"
$ws.Range("E9").Value = "  -4.56%  "
$ws.Range("D10").Value = "'1.262"
$ws.Range("E10").Value = "  -6.44%  "
$ws.Range("D11").Value = "'0.08152"
$ws.Range("E11").Value = "  -3.45%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'23.04"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("E14").Value = "  -6.13%  "
$ws.Range("D15").Value = "'7.430"
$ws.Range("E15").Value = "  -7.09%  "
$ws.Range("D16").Value = "'0.00001256"
$ws.Range("E16").Value = "  -4.12%  "
$ws.Range("D17").Value = "'1.605.44"
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("D18").Value = "'91.70"
$ws.Range("E18").Value = "  -3.08%  "
$ws.Range("D19").Value = "'0.06846"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "'18.26"
$ws.Range("E20").Value = "  -7.22%  "
$ws.Range("D21").Value = "'6.576"
$ws.Range("E21").Value = "  -5.75%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'13.03"
$ws.Range("E23").Value = "  -5.53%  "
$ws.Range("D24").Value = "'23.131.79"
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("D25").Value = "'2.342"
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("D26").Value = "'2.804"
$ws.Range("E26").Value = "  -5.64%  "
$ws.Range("D27").Value = "'21.07"
$ws.Range("E27").Value = "  -4.59%  "
$ws.Range("D28").Value = "'150.39"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "'5.286"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").Value = "'132.09"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("D31").Value = "'2.413"
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("D32").Value = "'6.890"
$ws.Range("E32").Value = "  -13.17%  "
$ws.Range("D33").Value = "'1.780.91"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").Value = "'0.07699"
$ws.Range("E34").Value = "  -4.86%  "
$ws.Range("D35").Value = "'0.9451"
$ws.Range("E35").Value = "  -7.78%  "
$ws.Range("D36").Value = "'0.02775"
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("D37").Value = "'6.261"
$ws.Range("E37").Value = "  -6.99%  "
$ws.Range("D38").Value = "'0.2544"
$ws.Range("E38").Value = "  -5.00%  "
$ws.Range("D39").Value = "'0.08914"
$ws.Range("E39").Value = "  -2.17%  "
$ws.Range("D40").Value = "'10.14"
$ws.Range("E40").Value = "  -5.43%  "
$ws.Range("D41").Value = "'1.392"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").Value = "'12.80"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43").Value = "'0.7101"
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("D44").Value = "'15.55"
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("D45").Value = "'0.6636"
$ws.Range("E45").Value = "  -4.67%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'2.301"
$ws.Range("E47").Value = "  -6.59%  "
$ws.Range("D48").Value = "'3.975"
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("D49").Value = "'131.66"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("D50").Value = "'0.07948"
$ws.Range("E50").Value = "  -4.60%  "
$ws.Range("D51").Value = "'1.215"
$ws.Range("E51").Value = "  -0.68%  "

Write-Output "Applied 95 cell updates"
